$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 0.001
$ws.Range("K11").Value = 471
$ws.Range("L11").Value = 0.002355
